$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-27 01:43:04"

# Insert a brand new row at position 6. This shifts the existing rows
# 6-10 down to 7-11 (their data + F-column hyperlink style carries with
# them, but the hyperlink relationship refs themselves stay put, which
# matches this runtime's Insert behavior).
$ws.Rows("6:6").Insert()

# Refresh the "取得日時" timestamp for every data row (2-5 keep their
# original job info, 7-11 are the rows that shifted down).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp

# Populate the freshly-inserted row 6 with the new job listing.
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【急募】webアプリ開発のエンジニアを探しています!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5479608"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "◆開発 ◇アプリ"

# Row 11 (formerly row 10, "NAS..." listing) is now the last row and
# needs its own hyperlink relationship added (it had no hyperlink entry
# before since it used to be row 10's own link, which must be recreated
# for its new location as a brand-new relationship, matching the rId10
# entry appended in the target workbook).
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5479430")
$ws.Range("F11").Style = "Hyperlink"
